$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Needs")

# Add a new "uncategorized" entry to the Transaction Type list (row 8)
$ws.Range("E8").Value = "uncategorized"

# Remove the "Organization Type Must be one of the following" section
# (old rows 11-12) and shift the remaining sections up by two rows.

# Row 11: was "Organization Type Must be one of the following"
#   -> becomes "Blockchain Network Name must be one of the following"
$ws.Range("A11").Value = "Blockchain Network Name must be one of the following"

# Row 12: was exchange / non_fungible_token / ico/token
#   -> becomes bitcoin / ripple / ethereum
$ws.Range("B12").Value = "bitcoin"
$ws.Range("C12").Value = "ripple"
$ws.Range("D12").Value = "ethereum"

# Row 13: was "Blockchain Network Name must be one of the following"
#   -> becomes "US fee must be to two decimal places"
$ws.Range("A13").Value = "US fee must be to two decimal places"

# Row 14: was bitcoin / ripple / ethereum (B14:D14)
#   -> becomes just "Status must be a 0 (failure) or 1 (success)" in A14
$ws.Range("B14:D14").ClearContents()
$ws.Range("A14").Value = "Status must be a 0 (failure) or 1 (success)"

# Row 15: was "US fee must be to two decimal places"
#   -> becomes "Currency Names Must be"
$ws.Range("A15").Value = "Currency Names Must be"

# Row 16: was "Status must be a 0 (failure) or 1 (success)"
#   -> becomes "All Currency Names must be uppercase"
$ws.Range("A16").Value = "All Currency Names must be uppercase"

# Rows 17 and 18 (old "Currency Names Must be" / "All Currency Names must be
# uppercase") no longer exist - delete them entirely so the sheet shrinks
# from A1:I18 to A1:I16.
$ws.Range("A17:I18").EntireRow.Delete()

# Update the selected cell to match the saved view state
$ws.Range("E8").Select()
